$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "94.367.51"
$ws.Range("E2").Value = "  +2.40%  "
$ws.Range("D3").Value = "3.072.09"
$ws.Range("E3").Value = "  -0.50%  "
$ws.Range("E4").Value = "  -0.18%  "
$ws.Range("D5").Value = "'237.12"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.87%  "
$ws.Range("D6").Value = "'610.14"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.32%  "
$ws.Range("E7").Value = "  +1.20%  "
$ws.Range("D8").Value = "'0.376"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.59%  "
$ws.Range("D9").Value = "'0.999"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.07%  "
$ws.Range("E10").Value = "  +10.43%  "
$ws.Range("D11").Value = "3.068.19"
$ws.Range("E11").Value = "  -0.44%  "
$ws.Range("E12").Value = "  -1.66%  "
$ws.Range("D13").Value = "94.050.26"
$ws.Range("E13").Value = "  +2.02%  "
$ws.Range("D14").Value = "'0.0000240"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.35%  "
$ws.Range("D15").Value = "'33.73"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.45%  "
$ws.Range("E16").Value = "  -1.30%  "
$ws.Range("D17").Value = "3.642.62"
$ws.Range("E17").Value = "  -0.91%  "
$ws.Range("D18").Value = "3.061.86"
$ws.Range("E18").Value = "  -1.39%  "
$ws.Range("D19").Value = "'3.54"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -5.20%  "
$ws.Range("D20").Value = "'14.33"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.24%  "
$ws.Range("D21").Value = "'5.64"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.87%  "
$ws.Range("D22").Value = "'437.56"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.75%  "
$ws.Range("D23").Value = "'8.79"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -5.06%  "
$ws.Range("D24").Value = "'0.0000188"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.06%  "
$ws.Range("D25").Value = "'8.35"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Value = "'5.50"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.57%  "
$ws.Range("D27").Value = "'84.61"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.17%  "
$ws.Range("D28").Value = "'11.81"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.33%  "
$ws.Range("E29").Value = "  -0.53%  "
$ws.Range("E30").Value = "  +0.65%  "
$ws.Range("E31").Value = "  +9.79%  "
$ws.Range("D32").Value = "'0.177"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.90%  "
$ws.Range("E33").Value = "  -5.17%  "
$ws.Range("D34").Value = "'9.00"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.60%  "
$ws.Range("D35").Value = "'0.989"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.42%  "
$ws.Range("D36").Value = "'7.66"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.57%  "
$ws.Range("D37").Value = "'0.153"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.61%  "
$ws.Range("D38").Value = "'25.33"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.60%  "
$ws.Range("D39").Value = "'1.88"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.22%  "
$ws.Range("D40").Value = "'478.57"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.05%  "
$ws.Range("D41").Value = "'24.05"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.73%  "
$ws.Range("D42").Value = "'3.75"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.69%  "
$ws.Range("B43").Value = "PolygonEcosystemToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D43").Value = "'0.434"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.70%  "
$ws.Range("B44").Value = "Fetch.AI"
$ws.Range("C44").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D44").Value = "'1.26"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.48%  "
$ws.Range("E45").Value = "  -0.01%  "
$ws.Range("D46").Value = "'3.10"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -5.06%  "
$ws.Range("D47").Value = "'161.64"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.61%  "
$ws.Range("D48").Value = "'0.670"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.43%  "
$ws.Range("E49").Value = "  -2.37%  "
$ws.Range("E50").Value = "  -0.67%  "
$ws.Range("D51").Value = "'0.997"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.05%  "
